# Sesuaikan format upload dan seeder
# Update the "desa_id" reference codes on the EPIDEMI PENYAKIT sheet so the
# template lines up with the reseeded village codes (53.06.13.2001..2008
# instead of the old 53.06.13.2014..2021 range), and leave the A2:A9 column
# selected afterwards (that's the range that was just retyped).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPIDEMI PENYAKIT")

$newDesaIds = @(
    "53.06.13.2001",
    "53.06.13.2002",
    "53.06.13.2003",
    "53.06.13.2004",
    "53.06.13.2005",
    "53.06.13.2006",
    "53.06.13.2007",
    "53.06.13.2008"
)

for ($i = 0; $i -lt $newDesaIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newDesaIds[$i]
}

$ws.Range("A2:A9").Select() | Out-Null
